# ---------------------------------------------------------------------------
# Regression.xlsx maintenance edit:
#  - rename "Sheet3" -> "EditMultipleProject" and populate it with a
#    "Transform & Fetch" table (copied layout from MultipleProjectAndModule,
#    with Project/Module/SubModule values replaced by "<old> : <new>" text
#    where the new part is bold).
#  - add a new "Sheet1" right after it containing the bottom 4 rows of that
#    same table (a second fixture used by the same test case).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. Rename the old placeholder sheet -----------------------------------
$wsEdit = $wb.Worksheets.Item(8)
$wsEdit.Name = "EditMultipleProject"

# ---- 2. Add the new sheet right after it (Excel default-names it "Sheet1") -
$wsNew = $wb.Worksheets.Add($null, $wsEdit)

# ---- 3. Seed EditMultipleProject from MultipleProjectAndModule's layout ----
$wsSrc = $wb.Worksheets.Item("MultipleProjectAndModule")
$wsSrc.Range("A1:F6").Copy($wsEdit.Range("A1:F6"))

# ---- 4. Overwrite the Project / Module / Sub-Module columns with the -------
#         "<original> : <edited>" text, bold-facing the edited half.
function Set-EditedCell {
    param($range, [string]$prefix, [string]$edited)

    $range.Value = $prefix + $edited
    $start = $prefix.Length + 1
    $len = $edited.Length
    $range.Characters($start, $len).Font.Bold = $true
}

Set-EditedCell $wsEdit.Range("A2") "KagamiProjectName1 : " "KagamiProjectNamEdited1"
Set-EditedCell $wsEdit.Range("A3") "KagamiProjectName2 : " "KagamiProjectNamEdited2"
Set-EditedCell $wsEdit.Range("A4") "KagamiProjectName3 : " "KagamiProjectNamEdited3"
Set-EditedCell $wsEdit.Range("A5") "KagamiProjectName4 : " "KagamiProjectNamEdited4"
Set-EditedCell $wsEdit.Range("A6") "KagamiProjectName5 : " "KagamiProjectNamEdited5"

Set-EditedCell $wsEdit.Range("C2") "Module13 " ": EditedModule1"
Set-EditedCell $wsEdit.Range("E2") "SubModule11 : " "EditedSubModule1"

Set-EditedCell $wsEdit.Range("C4") "Module33 : " "EditedModule3"
Set-EditedCell $wsEdit.Range("C5") "Module42 : " "EditedModule4"

Set-EditedCell $wsEdit.Range("E3") "SubModule22 : " "EditedSubModule2"
Set-EditedCell $wsEdit.Range("E4") "SubModule33 : " "EditedSubModule3"
Set-EditedCell $wsEdit.Range("E5") "SubModule41 : " "EditedSubModule4"
Set-EditedCell $wsEdit.Range("E6") "SubModule53 : " "EditedSubModule5"

Set-EditedCell $wsEdit.Range("C6") "Module51 : " "EditedModule5"
Set-EditedCell $wsEdit.Range("C3") "Module21 " ": EditedModule2"

# ---- 5. Auto-size the columns for the new table -----------------------------
$wsEdit.Columns.Item("A:F").EntireColumn.AutoFit()

# ---- 6. Populate the new Sheet1 with the last 4 rows of the table ----------
$wsEdit.Range("A3:F6").Copy($wsNew.Range("A1:F4"))
$wsNew.Columns.Item("A:F").EntireColumn.AutoFit()

# ---- 7. Restore per-sheet selections --------------------------------------
$wsSrc.Range("D12").Select()
$wsEdit.Range("D8").Select()
$wsNew.Range("A1:F4").Select()

# EditMultipleProject is the sheet that should end up active/selected.
$wsEdit.Activate()
$wsEdit.Range("D8").Select()
